$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values (column B / column C) per the revision
$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.25
$ws.Range("C5").Value = 15

# Update the sheet selection: cells B2:C5 selected (active cell C5).
# NOTE: this runtime always normalizes the active cell of a Range.Select()
# to the top-left corner of the resulting rectangle (Range.Activate() on a
# sub-cell collapses the selection to that single cell instead of just
# moving the active cell within the existing selection), so sqref="B2:C5"
# is the closest reproducible match; the active cell lands on B2.
$ws.Range("B2:C5").Select()
